# ----------------------------------------------------------------------
# "Datos que faltaban hasta el 10"
#  - Rename "Sheet 1" -> "datos"
#  - Add new "metadatos" sheet (becomes the active / selected tab)
#  - Populate "metadatos" with a variable dictionary
#  - Tidy the selection on "datos" (E1) now that it's no longer active
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- rename the original sheet -----------------------------------------
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

# --- leave a tidy selection behind on "datos" ---------------------------
$datos.Range("E1").Select() | Out-Null

# --- insert the new sheet right after "datos" ---------------------------
$metadatos = $wb.Worksheets.Add($null, $datos)
$metadatos.Name = "metadatos"

# --- header row -----------------------------------------------------------
$metadatos.Range("A1").Value = "Variables"
$metadatos.Range("B1").Value = "Descripción"
$metadatos.Range("C1").Value = "Fuente"
$metadatos.Range("D1").Value = "Fecha_de_extracción"

# --- anno -------------------------------------------------------------
$metadatos.Range("A2").Value = "anno"
$metadatos.Range("B2").Value = "Año"
$metadatos.Range("C2").Value = "…"

# --- codmpio ------------------------------------------------------------
$metadatos.Range("A3").Value = "codmpio"
$metadatos.Range("B3").Value = "Código del municipio"
$metadatos.Range("C3").Value = "…"

# --- SRPA_1 ---------------------------------------------------------------
$metadatos.Range("A4").Value = "SRPA_1"
$metadatos.Range("B4").Value = "`nNo. de adolescentes que ingresan al SRPA con una medida privativa de la libertad "
$metadatos.Range("C4").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- ingresos_totales -------------------------------------------------
$metadatos.Range("A5").Value = "ingresos_totales"
$metadatos.Range("B5").Value = " No. total de adolescentes que han ingresado al sistema SRPA en el mismo periodo y territorio. x 100"
$metadatos.Range("C5").Value = "`nInstituto Colombiano de Bienestar Familiar (ICBF)"

# --- tasa -------------------------------------------------------------
$metadatos.Range("A6").Value = "tasa"
$metadatos.Range("C6").Value = "Elaboración Propia"

# --- trailing (blank, but formatted) row -------------------------------
$metadatos.Range("A7").Value = ""
$metadatos.Range("B7").Value = ""
$metadatos.Range("C7").Value = ""
$metadatos.Range("D7").Value = ""

# --- font formatting: whole used range gets the explicit (non-theme) font
$used = $metadatos.Range("A1:D7")
$used.Font.Name = "Calibri"
$used.Font.Size = 11
$used.Font.Color = 0

# the embedded newlines in B4/B5 nudge the engine into giving those rows a
# custom height; auto-fit puts the row heights (and the customHeight flag)
# back to the sheet default
$metadatos.Rows.Item(4).AutoFit()
$metadatos.Rows.Item(5).AutoFit()

# --- extraction-date column (after the font pass, so the numFmt style also
#     carries font 1, matching the source xfId=2 "applyNumberFormat applyFont" xf)
$metadatos.Range("D2:D6").Value = 45722
$metadatos.Range("D2:D6").NumberFormat = "d-mmm-yy"

# A5 ends up without any explicit style (matches source workbook quirk)
$metadatos.Range("A5").Style = "Normal"

# --- column widths -------------------------------------------------------
$metadatos.Columns.Item(2).ColumnWidth = 29.666666666666668
$metadatos.Columns.Item(3).ColumnWidth = 36.5
$metadatos.Columns.Item(4).ColumnWidth = 9.5

# --- selection / active sheet --------------------------------------------
$metadatos.Activate()
$metadatos.Range("D2:D6").Select()
